# Insert a new row at position 11 (shifts old rows 11..117 down to 12..118)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("11:11").Insert()

# The newly inserted row 11 is blank. Populate it by duplicating the row
# that is now at row 12 (the former row 11), then overwrite the cells that
# differ for this new data point (date, min/max/avg price, and $/Kg).
$srcRow = 12
$dstRow = 11

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($col in $cols) {
    $srcCell = $ws.Range("$col$srcRow")
    $dstCell = $ws.Range("$col$dstRow")
    $dstCell.Value2 = $srcCell.Value2
}

# Now set the values that differ for the new row
$ws.Range("D11").Value2 = 45163
$ws.Range("N11").Value2 = 38000
$ws.Range("O11").Value2 = 38000
$ws.Range("P11").Value2 = 38000
$ws.Range("S11").Value2 = 2111
